# "Generate Report for Archive"
#
# 1. The status text "Ready for handoff" becomes "In Translation" everywhere
#    it appears (Overview!E2:F2, Overview!E3:F3, zh-cn!C2:C3, de-de!C2:C3 -
#    they all share the same string).
# 2. The "Status" column is narrowed from ~17.22 chars to ~13.41 chars on the
#    Overview sheet (columns E and F) and on the zh-cn / de-de sheets
#    (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status values -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Narrow the Status columns -------------------------------------------------
# The target raw OOXML column width is 13.4101848602295 characters; the
# ColumnWidth setter here snaps to a 1/6-character grid, so 12.5 is the
# value that lands on the closest achievable grid point (13.33...).
$newColumnWidth = 12.5

$overview.Range("E1").ColumnWidth = $newColumnWidth
$overview.Range("F1").ColumnWidth = $newColumnWidth

$zhcn.Range("C1").ColumnWidth = $newColumnWidth

$dede.Range("C1").ColumnWidth = $newColumnWidth
